# Edit: wording tweak in the "Contribution from Team Member 2" (Tyler Thrash)
# and "Contribution from Team Member 3" (Prashant Shrestha) cells of the
# Project Journal / Signoff table.
#   - "puzzle_words in the database."   ->  "puzzle_words table of the database."
#   - the Word "_GoBack" last-edit bookmark moves from the Prashant Shrestha
#     paragraph to sit right after the newly typed "table of" text.
#
# Because the two affected paragraphs contain many runs with identical
# formatting, the most reliable way to reproduce the precise run layout
# (including the now-mid-run split of "...changed an a|ttribute..." that
# Word's editor produced as a side effect, and the bookmark relocation) is
# to rebuild each paragraph's content from scratch with InsertXML.

$d = $word.ActiveDocument

$p20 = $null
$p23 = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Also changed an attribute name in puzzle_words in the database.*") {
        $p20 = $p
    }
    elseif ($t -like "*Created*admin_edit_synonyms*page functionality and UI that changes with the page.*") {
        $p23 = $p
    }
}

if ($null -eq $p20) {
    throw "Could not locate Tyler Thrash's contribution paragraph"
}
if ($null -eq $p23) {
    throw "Could not locate Prashant Shrestha's contribution paragraph"
}

$frag20 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0089511F" w:rsidRDefault="00EC1E36" w:rsidP="006E073E"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Tyler Thrash: </w:t></w:r><w:r w:rsidR="004D13D4"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Fixed a bug where the add puzzle functionality would allow the user to add a puzzle with the same name. Added play button functionality, and contributed to delete button</w:t></w:r><w:r w:rsidR="00145CFC"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00145CFC"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>admin_edit_synonyms</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00145CFC"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> page</w:t></w:r><w:r w:rsidR="004D13D4"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> functionality. Refactored code for better code reuse. Also changed an a</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">ttribute name in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004D13D4"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>puzzle_words</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="004D13D4"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> table of</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> the database.</w:t></w:r><w:r w:rsidR="00145CFC"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$frag23 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="000B17D1" w:rsidRDefault="000B17D1" w:rsidP="000B17D1"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Prashant Shrestha:</w:t></w:r><w:r w:rsidR="006E073E"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00B814BB"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Changed some minor UI according to feedback from the professor (changing button text, error text colors), Created </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="00B814BB"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>admin_edit_synonyms</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="00B814BB"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> page functionality and UI that changes with the page. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$p20.Range.InsertXML($frag20)
$p23.Range.InsertXML($frag23)

Write-Output "done"
